# Generate Report for Handoff
# Adds a new tracked file (a00a7228-422b-48f1-b114-67c1f80c027f.md) as row 9
# to the Overview, zh-cn and de-de tables/sheets.

$wb = $excel.ActiveWorkbook

$fileId   = "a00a7228-422b-48f1-b114-67c1f80c027f"
$fileName = "$fileId.md"
$pathName = "e2e\$fileId.md"
$hoHash   = "73f25297e0291b637cf1e62f058f73fb1aaa5130"
$hoUrl    = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/$hoHash/e2e/$fileName"

# ---------------------------------------------------------------------------
# Sheet "Overview" (sheet1) -> row 9
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.ListRows.Add() | Out-Null

$wsOverview.Range("A9").Value = $fileName
$wsOverview.Range("C9").Value = ".md"
$wsOverview.Range("D9").Value = ""
$wsOverview.Range("E9").Value = "Ready for handoff"
$wsOverview.Range("F9").Value = "Ready for handoff"
$wsOverview.Range("G9").Value = "2016-10-19 16:47:58"
$wsOverview.Range("G9").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsOverview.Hyperlinks.Add($wsOverview.Range("B9"), $hoUrl, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, $pathName) | Out-Null

# ---------------------------------------------------------------------------
# Sheet "zh-cn" (sheet2) -> row 9
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$loZhCn = $wsZhCn.ListObjects.Item(1)
$loZhCn.ListRows.Add() | Out-Null

$wsZhCn.Range("B9").Value = ".md"
$wsZhCn.Range("C9").Value = "Ready for handoff"
$wsZhCn.Range("D9").Value = "e2e"
$wsZhCn.Range("E9").Value = "ht"
$wsZhCn.Range("F9").Value = "False"
$wsZhCn.Range("G9").Value = "$fileId.$hoHash.zh-cn.xlf"
$wsZhCn.Range("H9").Value = "2016-10-19 16:47:47"
$wsZhCn.Range("H9").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Range("I9").Value = ""
$wsZhCn.Range("J9").Value = ""
$wsZhCn.Range("K9").Value = "0001-01-01 00:00:00"
$wsZhCn.Range("K9").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Range("L9").Value = ""
$wsZhCn.Range("M9").Value = "True"
$wsZhCn.Range("N9").Value = ""
$wsZhCn.Range("O9").Value = "False"
$wsZhCn.Range("P9").Value = ""

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A9"), $hoUrl, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, $fileName) | Out-Null

# ---------------------------------------------------------------------------
# Sheet "de-de" (sheet3) -> row 9
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$loDeDe = $wsDeDe.ListObjects.Item(1)
$loDeDe.ListRows.Add() | Out-Null

$wsDeDe.Range("B9").Value = ".md"
$wsDeDe.Range("C9").Value = "Ready for handoff"
$wsDeDe.Range("D9").Value = "e2e"
$wsDeDe.Range("E9").Value = "ht"
$wsDeDe.Range("F9").Value = "False"
$wsDeDe.Range("G9").Value = "$fileId.$hoHash.de-de.xlf"
$wsDeDe.Range("H9").Value = "2016-10-19 16:47:58"
$wsDeDe.Range("H9").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDeDe.Range("I9").Value = ""
$wsDeDe.Range("J9").Value = ""
$wsDeDe.Range("K9").Value = "0001-01-01 00:00:00"
$wsDeDe.Range("K9").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDeDe.Range("L9").Value = ""
$wsDeDe.Range("M9").Value = "True"
$wsDeDe.Range("N9").Value = ""
$wsDeDe.Range("O9").Value = "False"
$wsDeDe.Range("P9").Value = ""

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A9"), $hoUrl, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, $fileName) | Out-Null

Write-Output "Handoff report row added for $fileName"
